# Fruta / hortaliza, semanal
# Insert 4 new weekly rows (new Lane Late / Navel Late quotes) above the
# existing data block, pushing the previous rows 222-252 down to 226-256.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 222 (shifts old 222:252 -> 226:256).
$ws.Range("A222:A225").EntireRow.Insert()

$newRows = @(
    @{ Row = 222; K = "Lane Late";  L = "Primera"; M = 240; N = 6000; O = 6500; P = 6250; S = 417 },
    @{ Row = 223; K = "Lane Late";  L = "Segunda"; M = 200; N = 5000; O = 5500; P = 5250; S = 350 },
    @{ Row = 224; K = "Navel Late"; L = "Primera"; M = 240; N = 6000; O = 6500; P = 6250; S = 417 },
    @{ Row = 225; K = "Navel Late"; L = "Segunda"; M = 240; N = 5000; O = 5500; P = 5250; S = 350 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 7
    $ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($row, 3).Value = "Ñuble"
    $ws.Cells.Item($row, 4).Value = 44474
    $ws.Cells.Item($row, 5).Value = 16
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102005
    $ws.Cells.Item($row, 10).Value = "Naranja"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/bandeja 15 kilos granel"
    $ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 15
}
